# Apply weekly fruit/vegetable price update: the rows (2-35) are
# re-shuffled so that each record (Fecha, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg) moves to a new
# row position, per the source->destination mapping below.
# (destinationRow -> sourceRow, i.e. after[dest] = before[src])

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 34
    3  = 33
    4  = 27
    5  = 24
    6  = 13
    7  = 28
    8  = 9
    9  = 17
    10 = 8
    11 = 10
    12 = 12
    13 = 20
    14 = 21
    15 = 14
    16 = 25
    17 = 30
    18 = 15
    19 = 2
    20 = 6
    21 = 11
    22 = 35
    23 = 23
    24 = 3
    25 = 29
    26 = 5
    27 = 22
    28 = 4
    29 = 26
    30 = 31
    31 = 16
    32 = 18
    33 = 19
    34 = 32
    35 = 7
}

# Columns that carry data which gets shuffled between rows.
$cols = @("D", "L", "M", "N", "O", "P", "S")

# Snapshot the original values of the affected columns for every row
# before any writes happen (since several rows read from one another).
$original = @{}
for ($r = 2; $r -le 35; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $original[$r] = $rowVals
}

# Write the shuffled values into their new destination rows.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $original[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
